$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8972892761230469
$ws.Range("B1").Value = 1.364560484886169
$ws.Range("D1").Value = 1.680899262428284
$ws.Range("E1").Value = 1.097819924354553
